$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 4765.615  # H8
$ws.Cells.Item(8, 9).Value = 4495.8184  # I8
$ws.Cells.Item(8, 10).Value = 6249.5  # J8
$ws.Cells.Item(8, 11).Value = 13487.4552  # K8
$ws.Cells.Item(8, 12).Value = 18748.5  # L8
$ws.Cells.Item(8, 13).Value = -13348.4552  # M8
$ws.Cells.Item(8, 14).Value = -19026.5  # N8

$ws.Cells.Item(32, 8).Value = 5594.125  # H32
$ws.Cells.Item(32, 9).Value = 4633.6665  # I32
$ws.Cells.Item(32, 10).Value = 6170.4  # J32
$ws.Cells.Item(32, 11).Value = 4633.6665  # K32
$ws.Cells.Item(32, 12).Value = 6170.4  # L32
$ws.Cells.Item(32, 13).Value = -4307.6665  # M32
$ws.Cells.Item(32, 14).Value = -6822.4  # N32

$ws.Cells.Item(43, 8).Value = 13890.363  # H43
$ws.Cells.Item(43, 9).Value = 15616.5  # I43
$ws.Cells.Item(43, 10).Value = 11819  # J43
$ws.Cells.Item(43, 11).Value = 15616.5  # K43
$ws.Cells.Item(43, 12).Value = 11819  # L43
$ws.Cells.Item(43, 13).Value = -15547.5  # M43
$ws.Cells.Item(43, 14).Value = -11957  # N43

$ws.Cells.Item(51, 8).Value = 20350  # H51
$ws.Cells.Item(51, 9).Value = 10785.571  # I51
$ws.Cells.Item(51, 10).Value = 25500.076  # J51
$ws.Cells.Item(51, 11).Value = 10785.571  # K51
$ws.Cells.Item(51, 12).Value = 25500.076  # L51
$ws.Cells.Item(51, 13).Value = -10301.571  # M51
$ws.Cells.Item(51, 14).Value = -26468.076  # N51

$ws.Cells.Item(137, 8).Value = 20837550  # H137
$ws.Cells.Item(137, 9).Value = 125000250  # I137
$ws.Cells.Item(137, 11).Value = 375000750  # K137
$ws.Cells.Item(137, 13).Value = -374998200  # M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 7150.2607  # H61
$ws.Cells.Item(61, 9).Value = 4448  # I61
$ws.Cells.Item(61, 11).Value = 4448  # K61
$ws.Cells.Item(61, 13).Value = -4236  # M61

$ws.Cells.Item(74, 8).Value = 562835.9  # H74
$ws.Cells.Item(74, 9).Value = 913822.6  # I74
$ws.Cells.Item(74, 11).Value = 913822.6  # K74
$ws.Cells.Item(74, 13).Value = -912948.6  # M74

$ws.Cells.Item(77, 8).Value = 562835.9  # H77
$ws.Cells.Item(77, 9).Value = 913822.6  # I77
$ws.Cells.Item(77, 11).Value = 4569113  # K77
$ws.Cells.Item(77, 13).Value = -4564745  # M77

$ws.Cells.Item(122, 8).Value = 3457.0667  # H122
$ws.Cells.Item(122, 9).Value = 3527.3845  # I122
$ws.Cells.Item(122, 11).Value = 10582.1535  # K122
$ws.Cells.Item(122, 13).Value = -8132.1535  # M122

$ws.Cells.Item(136, 8).Value = 7150.2607  # H136
$ws.Cells.Item(136, 9).Value = 4448  # I136
$ws.Cells.Item(136, 11).Value = 13344  # K136
$ws.Cells.Item(136, 13).Value = -10794  # M136

$ws.Cells.Item(140, 8).Value = 75312.664  # H140
$ws.Cells.Item(140, 9).Value = 35000  # I140
$ws.Cells.Item(140, 10).Value = 83375.2  # J140
$ws.Cells.Item(140, 11).Value = 35000  # K140
$ws.Cells.Item(140, 12).Value = 83375.2  # L140
$ws.Cells.Item(140, 13).Value = -29820  # M140
$ws.Cells.Item(140, 14).Value = -93735.2  # N140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 28747  # H82
$ws.Cells.Item(82, 9).Value = 15710.929  # I82
$ws.Cells.Item(82, 11).Value = 15710.929  # K82
$ws.Cells.Item(82, 13).Value = -15327.929  # M82

$ws.Cells.Item(85, 8).Value = 28747  # H85
$ws.Cells.Item(85, 9).Value = 15710.929  # I85
$ws.Cells.Item(85, 11).Value = 15710.929  # K85
$ws.Cells.Item(85, 13).Value = -14384.929  # M85

$ws.Cells.Item(94, 8).Value = 4333.3335  # H94
$ws.Cells.Item(94, 9).Value = 4333.3335  # I94
$ws.Cells.Item(94, 11).Value = 4333.3335  # K94
$ws.Cells.Item(94, 13).Value = -3882.3335  # M94

$ws.Cells.Item(99, 8).Value = 3135.3333  # H99
$ws.Cells.Item(99, 9).Value = 2381.25  # I99
$ws.Cells.Item(99, 11).Value = 2381.25  # K99
$ws.Cells.Item(99, 13).Value = -883.25  # M99

$ws.Cells.Item(105, 8).Value = 90934270  # H105
$ws.Cells.Item(105, 9).Value = 111140860  # I105
$ws.Cells.Item(105, 11).Value = 111140860  # K105
$ws.Cells.Item(105, 13).Value = -111139113  # M105

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 41671090  # H31
$ws.Cells.Item(31, 9).Value = 100000920  # I31
$ws.Cells.Item(31, 10).Value = 6922.4287  # J31
$ws.Cells.Item(31, 11).Value = 100000920  # K31
$ws.Cells.Item(31, 12).Value = 6922.4287  # L31
$ws.Cells.Item(31, 13).Value = -100000625  # M31
$ws.Cells.Item(31, 14).Value = -7512.4287  # N31

$ws.Cells.Item(34, 8).Value = 41671090  # H34
$ws.Cells.Item(34, 9).Value = 100000920  # I34
$ws.Cells.Item(34, 10).Value = 6922.4287  # J34
$ws.Cells.Item(34, 11).Value = 100000920  # K34
$ws.Cells.Item(34, 12).Value = 6922.4287  # L34
$ws.Cells.Item(34, 13).Value = -100000718  # M34
$ws.Cells.Item(34, 14).Value = -7326.4287  # N34

$ws.Cells.Item(92, 8).Value = 50000  # H92
$ws.Cells.Item(92, 10).Value = 50000  # J92
$ws.Cells.Item(92, 12).Value = 50000  # L92
$ws.Cells.Item(92, 14).Value = -54992  # N92

$ws.Cells.Item(99, 8).Value = 1111  # H99
$ws.Cells.Item(99, 10).Value = 1111  # J99
$ws.Cells.Item(99, 12).Value = 1111  # L99
$ws.Cells.Item(99, 14).Value = -4107  # N99

$ws.Cells.Item(122, 8).Value = 112087.78  # H122
$ws.Cells.Item(122, 9).Value = 112087.78  # I122
$ws.Cells.Item(122, 11).Value = 336263.34  # K122
$ws.Cells.Item(122, 13).Value = -333813.34  # M122

$ws.Cells.Item(126, 8).Value = 1111  # H126
$ws.Cells.Item(126, 10).Value = 1111  # J126
$ws.Cells.Item(126, 12).Value = 3333  # L126
$ws.Cells.Item(126, 14).Value = -8273  # N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 787.7692  # H14
$ws.Cells.Item(14, 9).Value = 787.7692  # I14
$ws.Cells.Item(14, 11).Value = 2363.3076  # K14
$ws.Cells.Item(14, 13).Value = -2190.3076  # M14

$ws.Cells.Item(113, 8).Value = 3500  # H113
$ws.Cells.Item(113, 9).Value = 2000  # I113
$ws.Cells.Item(113, 11).Value = 6000  # K113
$ws.Cells.Item(113, 13).Value = -3830  # M113

$ws.Cells.Item(121, 8).Value = 15874219  # H121
$ws.Cells.Item(121, 10).Value = 25642572  # J121
$ws.Cells.Item(121, 12).Value = 76927716  # L121
$ws.Cells.Item(121, 14).Value = -76930336  # N121

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 14957.8  # H48
$ws.Cells.Item(48, 10).Value = 14957.8  # J48
$ws.Cells.Item(48, 12).Value = 14957.8  # L48
$ws.Cells.Item(48, 14).Value = -15927.8  # N48

$ws.Cells.Item(113, 8).Value = 10573.1  # H113
$ws.Cells.Item(113, 9).Value = 3144.8235  # I113
$ws.Cells.Item(113, 10).Value = 52666.668  # J113
$ws.Cells.Item(113, 11).Value = 3144.8235  # K113
$ws.Cells.Item(113, 12).Value = 52666.668  # L113
$ws.Cells.Item(113, 13).Value = -974.8235  # M113
$ws.Cells.Item(113, 14).Value = -57006.668  # N113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 127067.5  # H7
$ws.Cells.Item(7, 9).Value = 127067.5  # I7
$ws.Cells.Item(7, 11).Value = 127067.5  # K7
$ws.Cells.Item(7, 13).Value = -126955.5  # M7

$ws.Cells.Item(107, 8).Value = 2250.5  # H107
$ws.Cells.Item(107, 9).Value = 2250.5  # I107
$ws.Cells.Item(107, 11).Value = 2250.5  # K107
$ws.Cells.Item(107, 13).Value = -330.5  # M107

$ws.Cells.Item(122, 8).Value = 5103.391  # H122
$ws.Cells.Item(122, 9).Value = 4606.067  # I122
$ws.Cells.Item(122, 10).Value = 6035.875  # J122
$ws.Cells.Item(122, 11).Value = 13818.201  # K122
$ws.Cells.Item(122, 12).Value = 18107.625  # L122
$ws.Cells.Item(122, 13).Value = -11368.201  # M122
$ws.Cells.Item(122, 14).Value = -23007.625  # N122

$ws.Cells.Item(126, 8).Value = 127067.5  # H126
$ws.Cells.Item(126, 9).Value = 127067.5  # I126
$ws.Cells.Item(126, 11).Value = 381202.5  # K126
$ws.Cells.Item(126, 13).Value = -378732.5  # M126

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 0  # H4
$ws.Cells.Item(4, 10).Value = 0  # J4
$ws.Cells.Item(4, 12).Value = 0  # L4
$ws.Cells.Item(4, 14).ClearContents()  # N4

$ws.Cells.Item(132, 8).Value = 4414.3335  # H132
$ws.Cells.Item(132, 9).Value = 2193.1853  # I132
$ws.Cells.Item(132, 11).Value = 6579.5559  # K132
$ws.Cells.Item(132, 13).Value = -4049.5559  # M132
